# Populate the header row for the vaccinated-students sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen column E (Vaccine_Dose) to fit its header, matching the authored layout.
$ws.Columns.Item(5).ColumnWidth = 13.7

# Leave the active selection where the author left it after entering the data.
$ws.Range("F3").Select() | Out-Null
